$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B5: value changed from 3 to 0 ("use new seprate str instead of \t")
$ws.Range("B5").Value = 0

# C14: value removed entirely, cell left blank but keeps its style/format
# ("support excel ommit field")
$ws.Range("C14").ClearContents()

# Update the selected/active cell to reflect the saved view state
$ws.Range("D10").Select() | Out-Null
